# Update "想去人数" (interested-count) figures on the "展览" and "全部类型"
# sheets to reflect the latest generated output (see commit message).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 5435
$ws1.Range("F8").Value = 895
$ws1.Range("F9").Value = 138
$ws1.Range("F10").Value = 2419
$ws1.Range("F11").Value = 79
$ws1.Range("F12").Value = 61
$ws1.Range("F13").Value = 2265
$ws1.Range("F14").Value = 117

# --- Sheet "全部类型" ---------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 5435
$ws4.Range("F10").Value = 895
$ws4.Range("F11").Value = 138
$ws4.Range("F12").Value = 2419
$ws4.Range("F13").Value = 79
$ws4.Range("F15").Value = 61
$ws4.Range("F16").Value = 2265
$ws4.Range("F17").Value = 117

$wb.Save()
